# Applies the Phase-2 deck revision:
#  - Slide 1: reposition the 3rd name column and add "Sirul " before "Velaga"
#  - Slide 2 (agenda): title autosizes to 36pt; "Structure" bullet -> "Interface"
#  - Slide 3 (Problem statement): title autosizes to 36pt
#  - Slide 4: title "Structure" -> "Interface", autosizes to 36pt
#  - Slide 5 (Demo): title autosizes to 36pt
#  - Slide 6 (Algorithms): title autosizes to 36pt; last bullet "- " -> "- Chromatic number"
#  - Slide 9 (Plan for phase 3): merge runs, autosize to 36pt
#  - Remove slide 7 (the "Suggest a color to an uncolored node" pseudocode slide)

$p = $ppt.ActivePresentation
$ppAutoSizeTextToFitShape = 2

# ---------------------------------------------------------------------------
# Slide 1 - title slide: move the 3rd names column and insert "Sirul "
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$namesShape = $s1.Shapes.Item(3)
$namesShape.Left = 611.4174215748031
$namesShape.Top = 297.3909648818898

$tr1 = $namesShape.TextFrame.TextRange
$full1 = $tr1.Text
$idx1 = $full1.IndexOf("Velaga") + 1
$tr1.Characters($idx1, 6).Text = "Sirul Velaga"

# ---------------------------------------------------------------------------
# Slide 2 - agenda
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1)
$title2.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title2.TextFrame.TextRange.Font.Size = 36

$body2 = $s2.Shapes.Item(2)
$tr2 = $body2.TextFrame.TextRange
$full2 = $tr2.Text
$idx2 = $full2.IndexOf("Structure") + 1
$tr2.Characters($idx2, 9).Text = "Interface"

# ---------------------------------------------------------------------------
# Slide 3 - Problem statement
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title3.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------------------
# Slide 4 - Structure -> Interface
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "Interface"
$title4.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title4.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------------------
# Slide 5 - Demo
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Item(1)
$title5.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title5.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------------------
# Slide 6 - Algorithms
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$title6.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title6.TextFrame.TextRange.Font.Size = 36

$body6 = $s6.Shapes.Item(2)
$tr6 = $body6.TextFrame.TextRange
$full6 = $tr6.Text
$lastDashIdx = $full6.LastIndexOf("- ") + 1
$tr6.Characters($lastDashIdx, 2).Text = "- Chromatic number"

# ---------------------------------------------------------------------------
# Slide 9 - Plan for phase 3 (becomes slide 8 once slide 7 is removed below)
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1)
$title9.TextFrame.TextRange.Text = "Plan for phase 3"
$title9.TextFrame.AutoSize = $ppAutoSizeTextToFitShape
$title9.TextFrame.TextRange.Font.Size = 36

# ---------------------------------------------------------------------------
# Remove the pseudo-code slide ("Suggest a color to an uncolored node")
# ---------------------------------------------------------------------------
$p.Slides.Item(7).Delete()
